# Apply updated cryptocurrency price/volume data per commit "Updated cryptos list"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.421.10"
$ws.Range("E2").Value = "  -1.13%  "
$ws.Range("D3").Value = "3.754.15"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.58"
$ws.Range("E5").Value = "  -0.88%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "169.36"
$ws.Range("E6").Value = "  -0.13%  "
$ws.Range("D7").Value = "3.753.83"
$ws.Range("E7").Value = "  -1.93%  "
$ws.Range("E8").Value = "  -0.09%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.524"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("E10").Value = "  +0.47%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.48"
$ws.Range("E11").Value = "  +0.20%  "
$ws.Range("E12").Value = "  -1.09%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000276"
$ws.Range("E13").Value = "  +4.19%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.46"
$ws.Range("E14").Value = "  -1.81%  "
$ws.Range("D15").Value = "4.382.00"
$ws.Range("E15").Value = "  -2.16%  "
$ws.Range("D16").Value = "3.751.74"
$ws.Range("E16").Value = "  -2.03%  "
$ws.Range("E17").Value = "  +0.59%  "
$ws.Range("D18").Value = "67.345.21"
$ws.Range("E18").Value = "  -1.41%  "
$ws.Range("E19").Value = "  -3.08%  "
$ws.Range("E20").Value = "  +0.78%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.51"
$ws.Range("E21").Value = "  -5.06%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "466.39"
$ws.Range("E22").Value = "  -0.45%  "
$ws.Range("E23").Value = "  -2.41%  "
$ws.Range("E24").Value = "  -8.34%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.62"
$ws.Range("E25").Value = "  +0.59%  "
$ws.Range("E26").Value = "  -1.41%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.15"
$ws.Range("E27").Value = "  -0.60%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.29"
$ws.Range("E28").Value = "  +2.64%  "
$ws.Range("E30").Value = "  -2.41%  "
$ws.Range("D31").Value = "3.902.41"
$ws.Range("E31").Value = "  -1.93%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.65"
$ws.Range("E32").Value = "  -0.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.45"
$ws.Range("E33").Value = "  -3.53%  "
$ws.Range("E34").Value = "  -3.82%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.13"
$ws.Range("E35").Value = "  -2.96%  "
$ws.Range("D36").Value = "3.715.79"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.79"
$ws.Range("E37").Value = "  +2.75%  "
$ws.Range("E38").Value = "  -1.52%  "
$ws.Range("E39").Value = "  -2.19%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.998"
$ws.Range("E40").Value = "  -2.20%  "
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("E42").Value = "  -0.20%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.312"
$ws.Range("E43").Value = "  -0.97%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.71"
$ws.Range("E45").Value = "  -0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.95"
$ws.Range("E46").Value = "  -1.86%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "45.83"
$ws.Range("E47").Value = "  -2.81%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "397.76"
$ws.Range("E48").Value = "  -4.74%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000270"
$ws.Range("E49").Value = "  -7.50%  "
$ws.Range("B50").Value = "VeChain"
$ws.Range("C50").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0353"
$ws.Range("E50").Value = "  -1.87%  "
$ws.Range("B51").Value = "Monero"
$ws.Range("C51").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "138.99"
$ws.Range("E51").Value = "  -2.02%  "
